$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Planner's cost" map for problems 12-16 (rows 13-17, column H)
$ws.Range("H13").Value = -1
$ws.Range("H14").Value = -1
$ws.Range("H15").Value = -1
$ws.Range("H16").Value = -1
$ws.Range("H17").Value = 429

# Register the hidden "quick analysis" chart-data defined names that Excel
# creates when a chart is built from columns G and H (two chart revisions,
# v1 and v2)
$n = $wb.Names.Add('_xlchart.v1.4', 'Sheet1!$G$1')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v1.5', 'Sheet1!$G$2:$G$17')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v1.6', 'Sheet1!$H$1')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v1.7', 'Sheet1!$H$2:$H$17')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v2.0', 'Sheet1!$G$1')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v2.1', 'Sheet1!$G$2:$G$17')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v2.2', 'Sheet1!$H$1')
$n.Visible = $false
$n = $wb.Names.Add('_xlchart.v2.3', 'Sheet1!$H$2:$H$17')
$n.Visible = $false

# Match the final active selection recorded in the saved workbook
[void]$ws.Range("H17").Select()
